$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B8: change from inline string "2" to a numeric value 2
$ws.Range("B8").Value = 2

# Add new row 9 data
$ws.Range("A9").Value = "Ruilin"

# B9 must stay as text "4" (not auto-converted to a number)
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "4"
$ws.Range("B9").Style = "Normal"

$ws.Range("C9").Value = "Thank, thoughtful feedback"
$ws.Range("D9").Value = "ACK"
$ws.Range("E9").Value = "OTH"
$ws.Range("F9").Value = "41c93df3-3a59-4ce4-b94b-f420b7540586"
$ws.Range("G9").Value = "SJ19eUg0-_annotated.xlsx"
$ws.Range("H9").Value = "Thank the reviewer for the thoughtful feedback."
